# Applies:
#  1) Colors the "Sorting by Owner..." bullet red (C9211E), matching its
#     sibling bullets in the "Show Artifact Owners" section.
#  2) Merges the two runs of the "Add an Unassign button..." bullet into one.
#  3) Colors six bullets in the "Repair & Condition Management" section red
#     (C9211E), matching the other sections, and splits the "Assignment
#     rule" bullet's text into two runs (same formatting on both).

$d = $word.ActiveDocument

# wdColor value for RGB C9211E (Word stores colors as 0xBBGGRR)
$RED = 1974729

function Set-ParagraphRed($paraTextPattern) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $paraTextPattern) {
            $p.Range.Font.Color = $RED
            return $true
        }
    }
    return $false
}

# 1) "Sorting by Owner should work (unowned last)."
Set-ParagraphRed "Sorting by Owner should work*" | Out-Null

# 2) Merge the two runs of the Unassign bullet into a single run.
$d.Content.Find.Execute(
    "Add an Unassign button (or row context menu). (unassign from wizard)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Add an Unassign button (or row context menu). (unassign from wizard)",
    2) | Out-Null

# 3) Repair & Condition Management bullets -> red
Set-ParagraphRed "Add a Repair action for artifacts:*" | Out-Null
Set-ParagraphRed "Users can increase an artifact*s condition by an entered amount.*" | Out-Null
Set-ParagraphRed "Condition is clamped to 0*100.*" | Out-Null
Set-ParagraphRed "Wear rule: each assignment reduces condition by 5*" | Out-Null
Set-ParagraphRed "Show condition in the UI and ensure tables refresh after changes.*" | Out-Null

# The "Assignment rule" bullet needs recoloring AND to be split into two
# runs: "Assignment rule: if condition < 10, assignment is blocked " and
# "with a clear message." (both red, identical formatting).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Assignment rule: if condition*") {
        # Color the whole paragraph first (this also colors the paragraph
        # mark's rPr, i.e. pPr/rPr, matching the other bullets).
        $p.Range.Font.Color = $RED

        $splitText = "Assignment rule: if condition < 10, assignment is blocked "
        $start = $p.Range.Start
        $splitPos = $start + $splitText.Length
        $end = $p.Range.End

        $r1 = $d.Range($start, $splitPos)
        $r2 = $d.Range($splitPos, $end)

        # Apply distinct colors first to force the engine to keep the runs
        # split apart, then equalize the colors to the final red.
        $r1.Font.Color = $RED
        $r2.Font.Color = 0
        $r1.Font.Color = $RED
        $r2.Font.Color = $RED
        break
    }
}
